$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$dataRange = $ws.Range("C2:G19")
$dataRange.NumberFormat = "@"

$rowData = @{
    2 = @("PSV", "1.0", "7.4", "93%", "70%", "3.73")
    3 = @("Feyenoord", "1.0", "7.6", "77%", "67%", "3.33")
    4 = @("Twente", "1.6", "6.6", "74%", "64%", "2.90")
    5 = @("AZ", "1.3", "6.2", "90%", "63%", "3.13")
    6 = @("Ajax", "1.2", "5.6", "97%", "80%", "3.97")
    7 = @("Nijmegen", "1.6", "4.3", "93%", "67%", "3.43")
    8 = @("Utrecht", "1.5", "5.5", "77%", "40%", "2.80")
    9 = @("G. A. Eagle", "1.3", "5.0", "83%", "47%", "2.77")
    10 = @("Sparta Rotterdam", "1.2", "5.4", "90%", "57%", "2.93")
    11 = @("Heerenveen", "1.2", "5.5", "97%", "70%", "3.53")
    12 = @("Fortuna Sittard", "1.7", "4.9", "73%", "60%", "2.87")
    13 = @("Almere City", "2.1", "4.7", "74%", "49%", "2.65")
    14 = @("Zwolle", "2.0", "3.7", "90%", "63%", "3.27")
    15 = @("Heracles Almelo", "2.0", "4.5", "90%", "70%", "3.40")
    16 = @("Excelsior", "1.5", "3.5", "93%", "77%", "3.60")
    17 = @("Waalwijk", "1.7", "3.9", "73%", "50%", "2.63")
    18 = @("FC Volendam", "1.7", "4.1", "87%", "67%", "3.47")
    19 = @("Vitesse", "1.6", "5.4", "83%", "57%", "3.00")
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    $ws.Cells.Item($row, 2).Value = $vals[0]
    $ws.Cells.Item($row, 3).Value = $vals[1]
    $ws.Cells.Item($row, 4).Value = $vals[2]
    $ws.Cells.Item($row, 5).Value = $vals[3]
    $ws.Cells.Item($row, 6).Value = $vals[4]
    $ws.Cells.Item($row, 7).Value = $vals[5]
}

$dataRange.Style = "Normal"
